# Add two new rows (14 and 15) of localisation data to the "Пучиена 2" sheet,
# mirroring the existing "file marker row" + "data row" pattern already used
# for rows 10/11 (SCRIPT/G01P03A/um1503.ssb).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- formatting -----------------------------------------------------------
# Row 14 is a "file name marker" row like row 10 (bordered/bold style, no
# B/C/D/E content). Row 15 is a regular data row like row 11.
$ws.Range("A10:E10").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A11:E11").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)   # xlPasteFormats

$ws.Application.CutCopyMode = $false

# Match the row heights used by the equivalent rows in the sheet (ht="43.2").
$ws.Rows.Item(14).RowHeight = 43.2
$ws.Rows.Item(15).RowHeight = 43.2

# --- values -----------------------------------------------------------
# Row 14: only the script filename marker cell is populated.
# Row 15: filename, line number, English source, translated (RU) string and
# the converted/encoded string.
# (Values are written in the same order the original author's tool used, so
# the shared-strings table comes out in the same order.)
# NOTE: these strings contain a LITERAL backslash-n ("\n", two characters -
# backslash then letter n), matching the source file's convention for
# in-game line breaks. This is not a PowerShell newline escape (which would
# be the backtick form `n) - plain "\n" in a double-quoted string is already
# literal in PowerShell, so no escaping is required.
$ws.Cells.Item(14, 1).Value2 = "SCRIPT/G01P03A/um2506.ssb"
$ws.Cells.Item(15, 3).Value2 = " The grand master of all things\nbad? Nope, never heard of him."
$ws.Cells.Item(15, 1).Value2 = "SCRIPT/G01P03A/us0107.ssb"
$ws.Cells.Item(15, 5).Value2 = " Ãñàîä íàòóåñ âòåãï òàíïãï\nðìïöïãï? Îå-à, ÿ îéœåãï ï îæí îå èîàý."
$ws.Cells.Item(15, 4).Value2 = " Гранд мастер всего самого\nплохого? Не-а, я ничего о нём не знаю."
$ws.Cells.Item(15, 2).Value2 = 76

# --- selection --------------------------------------------------------
# The author's last selection after the edit was E15.
$ws.Range("E15").Select()
